$d = $word.ActiveDocument
$d.Content.Find.Execute("Desarrollo_Problema", $false, $false, $false, $false, $false, $true, 1, $false, "Situacion", 2)
